$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell C1 "timesteps" and copy style from B1
$ws.Range("C1").Value = "timesteps"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats

# Update existing data row
$ws.Range("A2").Value = "adult#001"
$ws.Range("B2").Value = 0.09
$ws.Range("C2").Value = 2400
